# 守備成績.xlsx — re-upload edit
#
# The source diff is the result of the workbook being re-opened/re-saved
# (sheet renamed + a different sheet left active when it was saved).
# Cell data and number formats are unchanged.
#
# Observable / reproducible changes:
#   1. The first sheet "統一7-Eleven獅" is renamed to "統一7-ELEVEn獅"
#      (capitalisation of "Eleven" -> "ELEVEn").
#   2. The workbook's active/selected sheet moves from sheet "味全龍"
#      (previously tabSelected) to the first sheet, which becomes the
#      active tab.
#
# (Window geometry, the x15ac:absPath hint, the xr:revisionPtr GUID and
# the cosmetic x14ac:dyDescent/defaultRowHeight values in the diff are
# artifacts of the authoring machine/display and are not meaningful,
# user-driven workbook edits reachable through the Excel object model.)

$wb = $excel.ActiveWorkbook

# 1. Rename the first worksheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "統一7-ELEVEn獅"

# 2. Make it the active sheet/tab (this both sets the workbook's
#    activeTab and moves tabSelected from the previously-active sheet).
$ws1.Activate()
